$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 7
    3  = 4
    4  = 8
    5  = 5
    6  = 4
    7  = 2
    8  = 4
    9  = 8
    10 = 4
    11 = 7
    12 = 10
    13 = 2
    14 = 6
    15 = 7
    16 = 3
    17 = 1
    18 = 3
    19 = 2
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
